$p = $ppt.ActivePresentation
$p.Slides.Item(19).Delete()
